# Generate Report for Handoff
# - Update Status text "Handed back: in sync with en-US" -> "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Narrow the date/status columns (was sized for the long "Handed back..." text)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-09-05 11:27:43"
$wsDeDe.Range("H2").Value = "2016-09-05 11:27:43"
$wsZhCn.Range("H2").Value = "2016-09-05 11:27:38"

# --- Column widths: shrink the status/date columns ---
# (Target stored width is 17.2159881591797 "characters"; ColumnWidth is
# quantized to the screen pixel grid by Excel, same as real COM automation,
# so 16.3 is the input that lands on the nearest representable width.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
